# Insert a new data row at row 81 (pushing the existing row 81..181 down to
# 82..182) and populate the new row with its own values, matching the
# weekly price update for "Fruta, Macroferia Regional de Talca - Piña".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 81..181 down to 82..182, leaving a blank row 81 that inherits
# the formatting (date style, etc.) of the row above it.
$ws.Rows.Item(81).Insert()

# Fill the newly inserted row 81 with the new record.
$ws.Range("A81").Value = 5
$ws.Range("B81").Value = "Macroferia Regional de Talca"
$ws.Range("C81").Value = "Maule"
$ws.Range("D81").Value = 44539
$ws.Range("E81").Value = 7
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100108
$ws.Range("H81").Value = "Tropicales y subtropicales"
$ws.Range("I81").Value = 100108005
$ws.Range("J81").Value = "Piña"
$ws.Range("K81").Value = "Caramelo"
$ws.Range("L81").Value = "Segunda"
$ws.Range("M81").Value = 200
$ws.Range("N81").Value = 18000
$ws.Range("O81").Value = 18000
$ws.Range("P81").Value = 18000
$ws.Range("Q81").Value = "`$/caja 14 unidades"
$ws.Range("R81").Value = "Ecuador"
$ws.Range("S81").Value = 1286
$ws.Range("T81").Value = 14
